# Update the TPM-derived NATMI output with the new values.
# Row 4 was moved out of the ECs block; "ECs" is now a full sending-cluster
# group (rows 2-4), and rows 8-10 add the missing MuSCs sending-cluster block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Ephb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05800433333333333
$ws.Range("H2").Value = 0.174013
$ws.Range("I2").Value = 0.02087975181349295
$ws.Range("J2").Value = 0.02087975181349295
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3407069999999999
$ws.Range("N2").Value = 1.022121
$ws.Range("O2").Value = 0.1055965976712818
$ws.Range("P2").Value = 0.1055965976712818
$ws.Range("Q2").Value = 0.019762482397
$ws.Range("R2").Value = 0.177862341573
$ws.Range("S2").Value = 0.002204830751725631
$ws.Range("T2").Value = 0.002204830751725631

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Ephb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05800433333333333
$ws.Range("H3").Value = 0.174013
$ws.Range("I3").Value = 0.02087975181349295
$ws.Range("J3").Value = 0.02087975181349295
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.220310333333333
$ws.Range("N3").Value = 3.660931
$ws.Range("O3").Value = 0.3782153560188308
$ws.Range("P3").Value = 0.3782153560188308
$ws.Range("Q3").Value = 0.07078328734477776
$ws.Range("R3").Value = 0.637049586103
$ws.Range("S3").Value = 0.007897042765725065
$ws.Range("T3").Value = 0.007897042765725065

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Ephb6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05800433333333333
$ws.Range("H4").Value = 0.174013
$ws.Range("I4").Value = 0.02087975181349295
$ws.Range("J4").Value = 0.02087975181349295
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.665478666666666
$ws.Range("N4").Value = 4.996435999999999
$ws.Range("O4").Value = 0.5161880463098875
$ws.Range("P4").Value = 0.5161880463098875
$ws.Range("Q4").Value = 0.09660497974088887
$ws.Range("R4").Value = 0.8694448176679999
$ws.Range("S4").Value = 0.01077787829604226
$ws.Range("T4").Value = 0.01077787829604226

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Ephb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.666083666666667
$ws.Range("H5").Value = 4.998251
$ws.Range("I5").Value = 0.5997381826733804
$ws.Range("J5").Value = 0.5997381826733805
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3407069999999999
$ws.Range("N5").Value = 1.022121
$ws.Range("O5").Value = 0.1055965976712818
$ws.Range("P5").Value = 0.1055965976712818
$ws.Range("Q5").Value = 0.5676463678189998
$ws.Range("R5").Value = 5.108817310370999
$ws.Range("S5").Value = 0.06333031158386666
$ws.Range("T5").Value = 0.06333031158386666

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Ephb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.666083666666667
$ws.Range("H6").Value = 4.998251
$ws.Range("I6").Value = 0.5997381826733804
$ws.Range("J6").Value = 0.5997381826733805
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.220310333333333
$ws.Range("N6").Value = 3.660931
$ws.Range("O6").Value = 0.3782153560188308
$ws.Range("P6").Value = 0.3782153560188308
$ws.Range("Q6").Value = 2.033139114631222
$ws.Range("R6").Value = 18.298252031681
$ws.Range("S6").Value = 0.2268301902778991
$ws.Range("T6").Value = 0.2268301902778992

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Ephb6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.666083666666667
$ws.Range("H7").Value = 4.998251
$ws.Range("I7").Value = 0.5997381826733804
$ws.Range("J7").Value = 0.5997381826733805
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.665478666666666
$ws.Range("N7").Value = 4.996435999999999
$ws.Range("O7").Value = 0.5161880463098875
$ws.Range("P7").Value = 0.5161880463098875
$ws.Range("Q7").Value = 2.77482680371511
$ws.Range("R7").Value = 24.973441233436
$ws.Range("S7").Value = 0.3095776808116146
$ws.Range("T7").Value = 0.3095776808116147

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Ephb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.053930333333333
$ws.Range("H8").Value = 3.161791
$ws.Range("I8").Value = 0.3793820655131266
$ws.Range("J8").Value = 0.3793820655131266
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3407069999999999
$ws.Range("N8").Value = 1.022121
$ws.Range("O8").Value = 0.1055965976712818
$ws.Range("P8").Value = 0.1055965976712818
$ws.Range("Q8").Value = 0.3590814420789999
$ws.Range("R8").Value = 3.231732978710999
$ws.Range("S8").Value = 0.0400614553356895
$ws.Range("T8").Value = 0.04006145533568949

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Ephb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.053930333333333
$ws.Range("H9").Value = 3.161791
$ws.Range("I9").Value = 0.3793820655131266
$ws.Range("J9").Value = 0.3793820655131266
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.220310333333333
$ws.Range("N9").Value = 3.660931
$ws.Range("O9").Value = 0.3782153560188308
$ws.Range("P9").Value = 0.3782153560188308
$ws.Range("Q9").Value = 1.286122076380111
$ws.Range("R9").Value = 11.575098687421
$ws.Range("S9").Value = 0.1434881229752066
$ws.Range("T9").Value = 0.1434881229752066

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efna5"
$ws.Range("C10").Value = "Ephb6"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.053930333333333
$ws.Range("H10").Value = 3.161791
$ws.Range("I10").Value = 0.3793820655131266
$ws.Range("J10").Value = 0.3793820655131266
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.665478666666666
$ws.Range("N10").Value = 4.996435999999999
$ws.Range("O10").Value = 0.5161880463098875
$ws.Range("P10").Value = 0.5161880463098875
$ws.Range("Q10").Value = 1.755298486319555
$ws.Range("R10").Value = 15.797686376876
$ws.Range("S10").Value = 0.1958324872022305
$ws.Range("T10").Value = 0.1958324872022305
